$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update swapped match rows (team/odds data swapped between paired rows, id/date stay fixed) ---
# Row 12
$ws.Cells.Item(12,2).Value = 6533424
$ws.Cells.Item(12,6).Value = "Kolding IF"
$ws.Cells.Item(12,7).Value = "Esbjerg"
$ws.Cells.Item(12,8).Value = 1
$ws.Cells.Item(12,9).Value = 2
$ws.Cells.Item(12,10).Value = "A"
$ws.Cells.Item(12,11).Value = 2.5
$ws.Cells.Item(12,12).Value = 3.75
$ws.Cells.Item(12,13).Value = 2.3
$ws.Cells.Item(12,14).Value = 2.25
$ws.Cells.Item(12,15).Value = 4
$ws.Cells.Item(12,16).Value = 2.5
$ws.Cells.Item(12,17).Value = 0
$ws.Cells.Item(12,18).Value = 1.775
$ws.Cells.Item(12,19).Value = 2.025
$ws.Cells.Item(12,20).Value = 2.75
$ws.Cells.Item(12,21).Value = 1.9
$ws.Cells.Item(12,22).Value = 1.9
$ws.Cells.Item(12,23).Value = -1
$ws.Cells.Item(12,24).Value = -1
$ws.Cells.Item(12,25).Value = 1.5
$ws.Cells.Item(12,26).Value = -1
$ws.Cells.Item(12,27).Value = 1.025
$ws.Cells.Item(12,28).Value = 0.45
$ws.Cells.Item(12,29).Value = -0.5

# Row 15
$ws.Cells.Item(15,2).Value = 6529284
$ws.Cells.Item(15,6).Value = "Frem"
$ws.Cells.Item(15,7).Value = "Brabrand"
$ws.Cells.Item(15,8).Value = 1
$ws.Cells.Item(15,9).Value = 1
$ws.Cells.Item(15,10).Value = "D"
$ws.Cells.Item(15,11).Value = 2.375
$ws.Cells.Item(15,12).Value = 3.8
$ws.Cells.Item(15,13).Value = 2.375
$ws.Cells.Item(15,14).Value = 2.4
$ws.Cells.Item(15,15).Value = 3.8
$ws.Cells.Item(15,16).Value = 2.3
$ws.Cells.Item(15,17).Value = 0
$ws.Cells.Item(15,18).Value = 1.975
$ws.Cells.Item(15,19).Value = 1.825
$ws.Cells.Item(15,20).Value = 3
$ws.Cells.Item(15,21).Value = 1.95
$ws.Cells.Item(15,22).Value = 1.85
$ws.Cells.Item(15,23).Value = -1
$ws.Cells.Item(15,24).Value = 2.8
$ws.Cells.Item(15,25).Value = -1
$ws.Cells.Item(15,26).Value = 0
$ws.Cells.Item(15,27).Value = -0
$ws.Cells.Item(15,28).Value = -1
$ws.Cells.Item(15,29).Value = 0.8500000000000001

# Row 17
$ws.Cells.Item(17,2).Value = 6858900
$ws.Cells.Item(17,6).Value = "Fremad Amager"
$ws.Cells.Item(17,7).Value = "Hellerup IK"
$ws.Cells.Item(17,8).Value = 2
$ws.Cells.Item(17,9).Value = 1
$ws.Cells.Item(17,10).Value = "H"
$ws.Cells.Item(17,11).Value = 1.833
$ws.Cells.Item(17,12).Value = 3.6
$ws.Cells.Item(17,13).Value = 3.4
$ws.Cells.Item(17,14).Value = 1.833
$ws.Cells.Item(17,15).Value = 3.6
$ws.Cells.Item(17,16).Value = 3.4
$ws.Cells.Item(17,17).Value = -0.5
$ws.Cells.Item(17,18).Value = 1.875
$ws.Cells.Item(17,19).Value = 1.925
$ws.Cells.Item(17,20).Value = 2.75
$ws.Cells.Item(17,21).Value = 1.825
$ws.Cells.Item(17,22).Value = 1.975
$ws.Cells.Item(17,23).Value = 0.833
$ws.Cells.Item(17,24).Value = -1
$ws.Cells.Item(17,25).Value = -1
$ws.Cells.Item(17,26).Value = 0.875
$ws.Cells.Item(17,27).Value = -1
$ws.Cells.Item(17,28).Value = 0.4125
$ws.Cells.Item(17,29).Value = -0.5

# Row 18
$ws.Cells.Item(18,2).Value = 6858901
$ws.Cells.Item(18,6).Value = "Skive"
$ws.Cells.Item(18,7).Value = "Esbjerg"
$ws.Cells.Item(18,8).Value = 1
$ws.Cells.Item(18,9).Value = 4
$ws.Cells.Item(18,10).Value = "A"
$ws.Cells.Item(18,11).Value = 3.75
$ws.Cells.Item(18,12).Value = 3.75
$ws.Cells.Item(18,13).Value = 1.727
$ws.Cells.Item(18,14).Value = 4.333
$ws.Cells.Item(18,15).Value = 4
$ws.Cells.Item(18,16).Value = 1.6
$ws.Cells.Item(18,17).Value = 0.75
$ws.Cells.Item(18,18).Value = 2
$ws.Cells.Item(18,19).Value = 1.8
$ws.Cells.Item(18,20).Value = 3
$ws.Cells.Item(18,21).Value = 1.95
$ws.Cells.Item(18,22).Value = 1.85
$ws.Cells.Item(18,23).Value = -1
$ws.Cells.Item(18,24).Value = -1
$ws.Cells.Item(18,25).Value = 0.6000000000000001
$ws.Cells.Item(18,26).Value = -1
$ws.Cells.Item(18,27).Value = 0.8
$ws.Cells.Item(18,28).Value = 0.95
$ws.Cells.Item(18,29).Value = -1

# Row 37
$ws.Cells.Item(37,2).Value = 6858920
$ws.Cells.Item(37,6).Value = "Middelfart"
$ws.Cells.Item(37,7).Value = "Fremad Amager"
$ws.Cells.Item(37,8).Value = 1
$ws.Cells.Item(37,9).Value = 0
$ws.Cells.Item(37,10).Value = "H"
$ws.Cells.Item(37,11).Value = 2.15
$ws.Cells.Item(37,12).Value = 3.6
$ws.Cells.Item(37,13).Value = 2.8
$ws.Cells.Item(37,14).Value = 1.85
$ws.Cells.Item(37,15).Value = 4
$ws.Cells.Item(37,16).Value = 3.3
$ws.Cells.Item(37,17).Value = -0.5
$ws.Cells.Item(37,18).Value = 1.875
$ws.Cells.Item(37,19).Value = 1.925
$ws.Cells.Item(37,20).Value = 3
$ws.Cells.Item(37,21).Value = 1.9
$ws.Cells.Item(37,22).Value = 1.9
$ws.Cells.Item(37,23).Value = 0.8500000000000001
$ws.Cells.Item(37,24).Value = -1
$ws.Cells.Item(37,25).Value = -1
$ws.Cells.Item(37,26).Value = 0.875
$ws.Cells.Item(37,27).Value = -1
$ws.Cells.Item(37,28).Value = -1
$ws.Cells.Item(37,29).Value = 0.8999999999999999

# Row 38
$ws.Cells.Item(38,2).Value = 6858919
$ws.Cells.Item(38,6).Value = "FC Roskilde"
$ws.Cells.Item(38,7).Value = "Skive"
$ws.Cells.Item(38,8).Value = 1
$ws.Cells.Item(38,9).Value = 1
$ws.Cells.Item(38,10).Value = "D"
$ws.Cells.Item(38,11).Value = 1.769
$ws.Cells.Item(38,12).Value = 3.75
$ws.Cells.Item(38,13).Value = 3.75
$ws.Cells.Item(38,14).Value = 2.1
$ws.Cells.Item(38,15).Value = 3.5
$ws.Cells.Item(38,16).Value = 2.9
$ws.Cells.Item(38,17).Value = -0.25
$ws.Cells.Item(38,18).Value = 1.9
$ws.Cells.Item(38,19).Value = 1.9
$ws.Cells.Item(38,20).Value = 2.75
$ws.Cells.Item(38,21).Value = 1.95
$ws.Cells.Item(38,22).Value = 1.85
$ws.Cells.Item(38,23).Value = -1
$ws.Cells.Item(38,24).Value = 2.5
$ws.Cells.Item(38,25).Value = -1
$ws.Cells.Item(38,26).Value = -0.5
$ws.Cells.Item(38,27).Value = 0.45
$ws.Cells.Item(38,28).Value = -1
$ws.Cells.Item(38,29).Value = 0.8500000000000001

# Row 41
$ws.Cells.Item(41,2).Value = 6858924
$ws.Cells.Item(41,6).Value = "Brabrand"
$ws.Cells.Item(41,7).Value = "Hellerup IK"
$ws.Cells.Item(41,8).Value = 1
$ws.Cells.Item(41,9).Value = 2
$ws.Cells.Item(41,10).Value = "A"
$ws.Cells.Item(41,11).Value = 2.15
$ws.Cells.Item(41,12).Value = 3.5
$ws.Cells.Item(41,13).Value = 2.9
$ws.Cells.Item(41,14).Value = 2.3
$ws.Cells.Item(41,15).Value = 3.4
$ws.Cells.Item(41,16).Value = 2.625
$ws.Cells.Item(41,17).Value = 0
$ws.Cells.Item(41,18).Value = 1.775
$ws.Cells.Item(41,19).Value = 2.025
$ws.Cells.Item(41,20).Value = 2.75
$ws.Cells.Item(41,21).Value = 1.925
$ws.Cells.Item(41,22).Value = 1.875
$ws.Cells.Item(41,23).Value = -1
$ws.Cells.Item(41,24).Value = -1
$ws.Cells.Item(41,25).Value = 1.625
$ws.Cells.Item(41,26).Value = -1
$ws.Cells.Item(41,27).Value = 1.025
$ws.Cells.Item(41,28).Value = 0.4625
$ws.Cells.Item(41,29).Value = -0.5

# Row 42
$ws.Cells.Item(42,2).Value = 6858923
$ws.Cells.Item(42,6).Value = "Skive"
$ws.Cells.Item(42,7).Value = "FA 2000"
$ws.Cells.Item(42,8).Value = 2
$ws.Cells.Item(42,9).Value = 2
$ws.Cells.Item(42,10).Value = "D"
$ws.Cells.Item(42,11).Value = 2.1
$ws.Cells.Item(42,12).Value = 3.5
$ws.Cells.Item(42,13).Value = 2.9
$ws.Cells.Item(42,14).Value = 2.55
$ws.Cells.Item(42,15).Value = 3.4
$ws.Cells.Item(42,16).Value = 2.375
$ws.Cells.Item(42,17).Value = 0
$ws.Cells.Item(42,18).Value = 1.975
$ws.Cells.Item(42,19).Value = 1.825
$ws.Cells.Item(42,20).Value = 2.5
$ws.Cells.Item(42,21).Value = 1.825
$ws.Cells.Item(42,22).Value = 1.975
$ws.Cells.Item(42,23).Value = -1
$ws.Cells.Item(42,24).Value = 2.4
$ws.Cells.Item(42,25).Value = -1
$ws.Cells.Item(42,26).Value = 0
$ws.Cells.Item(42,27).Value = -0
$ws.Cells.Item(42,28).Value = 0.825
$ws.Cells.Item(42,29).Value = -1

# Row 43
$ws.Cells.Item(43,2).Value = 6858926
$ws.Cells.Item(43,6).Value = "Nykobing"
$ws.Cells.Item(43,7).Value = "Middelfart"
$ws.Cells.Item(43,8).Value = 0
$ws.Cells.Item(43,9).Value = 0
$ws.Cells.Item(43,10).Value = "D"
$ws.Cells.Item(43,11).Value = 2.6
$ws.Cells.Item(43,12).Value = 3.5
$ws.Cells.Item(43,13).Value = 2.35
$ws.Cells.Item(43,14).Value = 2.6
$ws.Cells.Item(43,15).Value = 3.5
$ws.Cells.Item(43,16).Value = 2.35
$ws.Cells.Item(43,17).Value = 0
$ws.Cells.Item(43,18).Value = 2
$ws.Cells.Item(43,19).Value = 1.8
$ws.Cells.Item(43,20).Value = 3
$ws.Cells.Item(43,21).Value = 1.975
$ws.Cells.Item(43,22).Value = 1.825
$ws.Cells.Item(43,23).Value = -1
$ws.Cells.Item(43,24).Value = 2.5
$ws.Cells.Item(43,25).Value = -1
$ws.Cells.Item(43,26).Value = 0
$ws.Cells.Item(43,27).Value = -0
$ws.Cells.Item(43,28).Value = -1
$ws.Cells.Item(43,29).Value = 0.825

# Row 44
$ws.Cells.Item(44,2).Value = 6858927
$ws.Cells.Item(44,6).Value = "Thisted FC"
$ws.Cells.Item(44,7).Value = "Aarhus Fremad"
$ws.Cells.Item(44,8).Value = 1
$ws.Cells.Item(44,9).Value = 3
$ws.Cells.Item(44,10).Value = "A"
$ws.Cells.Item(44,11).Value = 4.333
$ws.Cells.Item(44,12).Value = 3.6
$ws.Cells.Item(44,13).Value = 1.7
$ws.Cells.Item(44,14).Value = 4.75
$ws.Cells.Item(44,15).Value = 3.6
$ws.Cells.Item(44,16).Value = 1.65
$ws.Cells.Item(44,17).Value = 0.75
$ws.Cells.Item(44,18).Value = 2
$ws.Cells.Item(44,19).Value = 1.8
$ws.Cells.Item(44,20).Value = 3
$ws.Cells.Item(44,21).Value = 1.925
$ws.Cells.Item(44,22).Value = 1.875
$ws.Cells.Item(44,23).Value = -1
$ws.Cells.Item(44,24).Value = -1
$ws.Cells.Item(44,25).Value = 0.6499999999999999
$ws.Cells.Item(44,26).Value = -1
$ws.Cells.Item(44,27).Value = 0.8
$ws.Cells.Item(44,28).Value = 0.925
$ws.Cells.Item(44,29).Value = -1

# Row 46
$ws.Cells.Item(46,2).Value = 6858929
$ws.Cells.Item(46,6).Value = "FC Roskilde"
$ws.Cells.Item(46,7).Value = "AB Copenhagen"
$ws.Cells.Item(46,8).Value = 2
$ws.Cells.Item(46,9).Value = 0
$ws.Cells.Item(46,10).Value = "H"
$ws.Cells.Item(46,11).Value = 2.45
$ws.Cells.Item(46,12).Value = 3.5
$ws.Cells.Item(46,13).Value = 2.45
$ws.Cells.Item(46,14).Value = 2.375
$ws.Cells.Item(46,15).Value = 3.4
$ws.Cells.Item(46,16).Value = 2.55
$ws.Cells.Item(46,17).Value = 0
$ws.Cells.Item(46,18).Value = 1.85
$ws.Cells.Item(46,19).Value = 1.95
$ws.Cells.Item(46,20).Value = 2.5
$ws.Cells.Item(46,21).Value = 1.875
$ws.Cells.Item(46,22).Value = 1.925
$ws.Cells.Item(46,23).Value = 1.375
$ws.Cells.Item(46,24).Value = -1
$ws.Cells.Item(46,25).Value = -1
$ws.Cells.Item(46,26).Value = 0.8500000000000001
$ws.Cells.Item(46,27).Value = -1
$ws.Cells.Item(46,28).Value = -1
$ws.Cells.Item(46,29).Value = 0.925

# Row 48
$ws.Cells.Item(48,2).Value = 6858930
$ws.Cells.Item(48,6).Value = "Thisted FC"
$ws.Cells.Item(48,7).Value = "Skive"
$ws.Cells.Item(48,8).Value = 0
$ws.Cells.Item(48,9).Value = 0
$ws.Cells.Item(48,10).Value = "D"
$ws.Cells.Item(48,11).Value = 2.4
$ws.Cells.Item(48,12).Value = 3.5
$ws.Cells.Item(48,13).Value = 2.4
$ws.Cells.Item(48,14).Value = 2.25
$ws.Cells.Item(48,15).Value = 3.5
$ws.Cells.Item(48,16).Value = 2.55
$ws.Cells.Item(48,17).Value = 0
$ws.Cells.Item(48,18).Value = 1.775
$ws.Cells.Item(48,19).Value = 2.025
$ws.Cells.Item(48,20).Value = 2.5
$ws.Cells.Item(48,21).Value = 1.85
$ws.Cells.Item(48,22).Value = 1.95
$ws.Cells.Item(48,23).Value = -1
$ws.Cells.Item(48,24).Value = 2.5
$ws.Cells.Item(48,25).Value = -1
$ws.Cells.Item(48,26).Value = 0
$ws.Cells.Item(48,27).Value = -0
$ws.Cells.Item(48,28).Value = -1
$ws.Cells.Item(48,29).Value = 0.95

# Row 50
$ws.Cells.Item(50,2).Value = 6858933
$ws.Cells.Item(50,6).Value = "Esbjerg"
$ws.Cells.Item(50,7).Value = "Nykobing"
$ws.Cells.Item(50,8).Value = 3
$ws.Cells.Item(50,9).Value = 4
$ws.Cells.Item(50,10).Value = "A"
$ws.Cells.Item(50,11).Value = 1.35
$ws.Cells.Item(50,12).Value = 5
$ws.Cells.Item(50,13).Value = 6.5
$ws.Cells.Item(50,14).Value = 1.333
$ws.Cells.Item(50,15).Value = 5.5
$ws.Cells.Item(50,16).Value = 6
$ws.Cells.Item(50,17).Value = -1.5
$ws.Cells.Item(50,18).Value = 1.9
$ws.Cells.Item(50,19).Value = 1.9
$ws.Cells.Item(50,20).Value = 3.25
$ws.Cells.Item(50,21).Value = 1.925
$ws.Cells.Item(50,22).Value = 1.875
$ws.Cells.Item(50,23).Value = -1
$ws.Cells.Item(50,24).Value = -1
$ws.Cells.Item(50,25).Value = 5
$ws.Cells.Item(50,26).Value = -1
$ws.Cells.Item(50,27).Value = 0.8999999999999999
$ws.Cells.Item(50,28).Value = 0.925
$ws.Cells.Item(50,29).Value = -1

# Row 51
$ws.Cells.Item(51,2).Value = 6858934
$ws.Cells.Item(51,6).Value = "Aarhus Fremad"
$ws.Cells.Item(51,7).Value = "Fremad Amager"
$ws.Cells.Item(51,8).Value = 3
$ws.Cells.Item(51,9).Value = 1
$ws.Cells.Item(51,10).Value = "H"
$ws.Cells.Item(51,11).Value = 1.55
$ws.Cells.Item(51,12).Value = 4
$ws.Cells.Item(51,13).Value = 5
$ws.Cells.Item(51,14).Value = 1.5
$ws.Cells.Item(51,15).Value = 4.2
$ws.Cells.Item(51,16).Value = 5.25
$ws.Cells.Item(51,17).Value = -1
$ws.Cells.Item(51,18).Value = 1.8
$ws.Cells.Item(51,19).Value = 2
$ws.Cells.Item(51,20).Value = 3
$ws.Cells.Item(51,21).Value = 1.9
$ws.Cells.Item(51,22).Value = 1.9
$ws.Cells.Item(51,23).Value = 0.5
$ws.Cells.Item(51,24).Value = -1
$ws.Cells.Item(51,25).Value = -1
$ws.Cells.Item(51,26).Value = 0.8
$ws.Cells.Item(51,27).Value = -1
$ws.Cells.Item(51,28).Value = 0.8999999999999999
$ws.Cells.Item(51,29).Value = -1

# Row 59
$ws.Cells.Item(59,2).Value = 6858942
$ws.Cells.Item(59,6).Value = "Middelfart"
$ws.Cells.Item(59,7).Value = "AB Copenhagen"
$ws.Cells.Item(59,8).Value = 1
$ws.Cells.Item(59,9).Value = 0
$ws.Cells.Item(59,10).Value = "H"
$ws.Cells.Item(59,11).Value = 2.1
$ws.Cells.Item(59,12).Value = 3.5
$ws.Cells.Item(59,13).Value = 3
$ws.Cells.Item(59,14).Value = 2.1
$ws.Cells.Item(59,15).Value = 3.5
$ws.Cells.Item(59,16).Value = 3
$ws.Cells.Item(59,17).Value = -0.25
$ws.Cells.Item(59,18).Value = 1.85
$ws.Cells.Item(59,19).Value = 1.95
$ws.Cells.Item(59,20).Value = 2.5
$ws.Cells.Item(59,21).Value = 1.85
$ws.Cells.Item(59,22).Value = 1.95
$ws.Cells.Item(59,23).Value = 1.1
$ws.Cells.Item(59,24).Value = -1
$ws.Cells.Item(59,25).Value = -1
$ws.Cells.Item(59,26).Value = 0.8500000000000001
$ws.Cells.Item(59,27).Value = -1
$ws.Cells.Item(59,28).Value = -1
$ws.Cells.Item(59,29).Value = 0.95

# Row 60
$ws.Cells.Item(60,2).Value = 6858941
$ws.Cells.Item(60,6).Value = "Aarhus Fremad"
$ws.Cells.Item(60,7).Value = "FA 2000"
$ws.Cells.Item(60,8).Value = 3
$ws.Cells.Item(60,9).Value = 0
$ws.Cells.Item(60,10).Value = "H"
$ws.Cells.Item(60,11).Value = 1.5
$ws.Cells.Item(60,12).Value = 4
$ws.Cells.Item(60,13).Value = 5.5
$ws.Cells.Item(60,14).Value = 1.444
$ws.Cells.Item(60,15).Value = 4.2
$ws.Cells.Item(60,16).Value = 6
$ws.Cells.Item(60,17).Value = -1
$ws.Cells.Item(60,18).Value = 1.725
$ws.Cells.Item(60,19).Value = 2.075
$ws.Cells.Item(60,20).Value = 3
$ws.Cells.Item(60,21).Value = 1.95
$ws.Cells.Item(60,22).Value = 1.85
$ws.Cells.Item(60,23).Value = 0.444
$ws.Cells.Item(60,24).Value = -1
$ws.Cells.Item(60,25).Value = -1
$ws.Cells.Item(60,26).Value = 0.7250000000000001
$ws.Cells.Item(60,27).Value = -1
$ws.Cells.Item(60,28).Value = 0
$ws.Cells.Item(60,29).Value = -0

# Row 82
$ws.Cells.Item(82,2).Value = 6859008
$ws.Cells.Item(82,6).Value = "Brabrand"
$ws.Cells.Item(82,7).Value = "AB Copenhagen"
$ws.Cells.Item(82,8).Value = 2
$ws.Cells.Item(82,9).Value = 2
$ws.Cells.Item(82,10).Value = "D"
$ws.Cells.Item(82,11).Value = 3.6
$ws.Cells.Item(82,12).Value = 3.6
$ws.Cells.Item(82,13).Value = 1.85
$ws.Cells.Item(82,14).Value = 3.4
$ws.Cells.Item(82,15).Value = 3.6
$ws.Cells.Item(82,16).Value = 1.909
$ws.Cells.Item(82,17).Value = 0.5
$ws.Cells.Item(82,18).Value = 1.825
$ws.Cells.Item(82,19).Value = 1.975
$ws.Cells.Item(82,20).Value = 2.5
$ws.Cells.Item(82,21).Value = 1.8
$ws.Cells.Item(82,22).Value = 2
$ws.Cells.Item(82,23).Value = -1
$ws.Cells.Item(82,24).Value = 2.6
$ws.Cells.Item(82,25).Value = -1
$ws.Cells.Item(82,26).Value = 0.825
$ws.Cells.Item(82,27).Value = -1
$ws.Cells.Item(82,28).Value = 0.8
$ws.Cells.Item(82,29).Value = -1

# Row 83
$ws.Cells.Item(83,2).Value = 6859007
$ws.Cells.Item(83,6).Value = "Skive"
$ws.Cells.Item(83,7).Value = "FC Roskilde"
$ws.Cells.Item(83,8).Value = 1
$ws.Cells.Item(83,9).Value = 2
$ws.Cells.Item(83,10).Value = "A"
$ws.Cells.Item(83,11).Value = 3.6
$ws.Cells.Item(83,12).Value = 3.4
$ws.Cells.Item(83,13).Value = 1.909
$ws.Cells.Item(83,14).Value = 3.2
$ws.Cells.Item(83,15).Value = 3.4
$ws.Cells.Item(83,16).Value = 2.05
$ws.Cells.Item(83,17).Value = 0.25
$ws.Cells.Item(83,18).Value = 2
$ws.Cells.Item(83,19).Value = 1.8
$ws.Cells.Item(83,20).Value = 2.75
$ws.Cells.Item(83,21).Value = 1.975
$ws.Cells.Item(83,22).Value = 1.825
$ws.Cells.Item(83,23).Value = -1
$ws.Cells.Item(83,24).Value = -1
$ws.Cells.Item(83,25).Value = 1.05
$ws.Cells.Item(83,26).Value = -1
$ws.Cells.Item(83,27).Value = 0.8
$ws.Cells.Item(83,28).Value = 0.4875
$ws.Cells.Item(83,29).Value = -0.5

# Row 84
$ws.Cells.Item(84,2).Value = 6859010
$ws.Cells.Item(84,6).Value = "Esbjerg"
$ws.Cells.Item(84,7).Value = "FA 2000"
$ws.Cells.Item(84,8).Value = 3
$ws.Cells.Item(84,9).Value = 1
$ws.Cells.Item(84,10).Value = "H"
$ws.Cells.Item(84,11).Value = 1.222
$ws.Cells.Item(84,12).Value = 6.5
$ws.Cells.Item(84,13).Value = 9
$ws.Cells.Item(84,14).Value = 1.2
$ws.Cells.Item(84,15).Value = 6.5
$ws.Cells.Item(84,16).Value = 10
$ws.Cells.Item(84,17).Value = -2
$ws.Cells.Item(84,18).Value = 1.9
$ws.Cells.Item(84,19).Value = 1.9
$ws.Cells.Item(84,20).Value = 3.5
$ws.Cells.Item(84,21).Value = 1.875
$ws.Cells.Item(84,22).Value = 1.925
$ws.Cells.Item(84,23).Value = 0.2
$ws.Cells.Item(84,24).Value = -1
$ws.Cells.Item(84,25).Value = -1
$ws.Cells.Item(84,26).Value = 0
$ws.Cells.Item(84,27).Value = -0
$ws.Cells.Item(84,28).Value = 0.875
$ws.Cells.Item(84,29).Value = -1

# Row 85
$ws.Cells.Item(85,2).Value = 6859011
$ws.Cells.Item(85,6).Value = "Middelfart"
$ws.Cells.Item(85,7).Value = "Nykobing"
$ws.Cells.Item(85,8).Value = 2
$ws.Cells.Item(85,9).Value = 2
$ws.Cells.Item(85,10).Value = "D"
$ws.Cells.Item(85,11).Value = 2
$ws.Cells.Item(85,12).Value = 3.5
$ws.Cells.Item(85,13).Value = 3.3
$ws.Cells.Item(85,14).Value = 2
$ws.Cells.Item(85,15).Value = 3.5
$ws.Cells.Item(85,16).Value = 3.3
$ws.Cells.Item(85,17).Value = -0.25
$ws.Cells.Item(85,18).Value = 1.75
$ws.Cells.Item(85,19).Value = 1.95
$ws.Cells.Item(85,20).Value = 2.5
$ws.Cells.Item(85,21).Value = 1.8
$ws.Cells.Item(85,22).Value = 2
$ws.Cells.Item(85,23).Value = -1
$ws.Cells.Item(85,24).Value = 2.5
$ws.Cells.Item(85,25).Value = -1
$ws.Cells.Item(85,26).Value = -0.5
$ws.Cells.Item(85,27).Value = 0.475
$ws.Cells.Item(85,28).Value = 0.8
$ws.Cells.Item(85,29).Value = -1

# --- Append new rows 112-117 ---
# Copy formatting for id (A) and date (E) columns from row 111 so new rows match existing style
# Row 112
$ws.Cells.Item(111,1).Copy() | Out-Null
$ws.Cells.Item(112,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(111,5).Copy() | Out-Null
$ws.Cells.Item(112,5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(112,1).Value = 110
$ws.Cells.Item(112,2).Value = 6859072
$ws.Cells.Item(112,3).Value = "Denmark Division 2"
$ws.Cells.Item(112,4).Value = "Denmark Division 2"
$ws.Cells.Item(112,5).Value = 45360.40625
$ws.Cells.Item(112,6).Value = "FA 2000"
$ws.Cells.Item(112,7).Value = "Brabrand"
$ws.Cells.Item(112,8).Value = 2
$ws.Cells.Item(112,9).Value = 0
$ws.Cells.Item(112,10).Value = "H"
$ws.Cells.Item(112,11).Value = 2.1
$ws.Cells.Item(112,12).Value = 3.5
$ws.Cells.Item(112,13).Value = 2.9
$ws.Cells.Item(112,14).Value = 1.95
$ws.Cells.Item(112,15).Value = 3.6
$ws.Cells.Item(112,16).Value = 3.2
$ws.Cells.Item(112,17).Value = -0.5
$ws.Cells.Item(112,18).Value = 2
$ws.Cells.Item(112,19).Value = 1.8
$ws.Cells.Item(112,20).Value = 2.5
$ws.Cells.Item(112,21).Value = 1.85
$ws.Cells.Item(112,22).Value = 1.95
$ws.Cells.Item(112,23).Value = 0.95
$ws.Cells.Item(112,24).Value = -1
$ws.Cells.Item(112,25).Value = -1
$ws.Cells.Item(112,26).Value = 1
$ws.Cells.Item(112,27).Value = -1
$ws.Cells.Item(112,28).Value = -1
$ws.Cells.Item(112,29).Value = 0.95

# Row 113
$ws.Cells.Item(111,1).Copy() | Out-Null
$ws.Cells.Item(113,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(111,5).Copy() | Out-Null
$ws.Cells.Item(113,5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(113,1).Value = 111
$ws.Cells.Item(113,2).Value = 6859074
$ws.Cells.Item(113,3).Value = "Denmark Division 2"
$ws.Cells.Item(113,4).Value = "Denmark Division 2"
$ws.Cells.Item(113,5).Value = 45360.41666666666
$ws.Cells.Item(113,6).Value = "FC Roskilde"
$ws.Cells.Item(113,7).Value = "Middelfart"
$ws.Cells.Item(113,8).Value = 3
$ws.Cells.Item(113,9).Value = 1
$ws.Cells.Item(113,10).Value = "H"
$ws.Cells.Item(113,11).Value = 2.2
$ws.Cells.Item(113,12).Value = 3.4
$ws.Cells.Item(113,13).Value = 2.8
$ws.Cells.Item(113,14).Value = 2.15
$ws.Cells.Item(113,15).Value = 3.4
$ws.Cells.Item(113,16).Value = 2.875
$ws.Cells.Item(113,17).Value = -0.25
$ws.Cells.Item(113,18).Value = 1.95
$ws.Cells.Item(113,19).Value = 1.85
$ws.Cells.Item(113,20).Value = 2.5
$ws.Cells.Item(113,21).Value = 1.975
$ws.Cells.Item(113,22).Value = 1.825
$ws.Cells.Item(113,23).Value = 1.15
$ws.Cells.Item(113,24).Value = -1
$ws.Cells.Item(113,25).Value = -1
$ws.Cells.Item(113,26).Value = 0.95
$ws.Cells.Item(113,27).Value = -1
$ws.Cells.Item(113,28).Value = 0.9750000000000001
$ws.Cells.Item(113,29).Value = -1

# Row 114
$ws.Cells.Item(111,1).Copy() | Out-Null
$ws.Cells.Item(114,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(111,5).Copy() | Out-Null
$ws.Cells.Item(114,5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(114,1).Value = 112
$ws.Cells.Item(114,2).Value = 6859073
$ws.Cells.Item(114,3).Value = "Denmark Division 2"
$ws.Cells.Item(114,4).Value = "Denmark Division 2"
$ws.Cells.Item(114,5).Value = 45360.41666666666
$ws.Cells.Item(114,6).Value = "Aarhus Fremad"
$ws.Cells.Item(114,7).Value = "Thisted FC"
$ws.Cells.Item(114,8).Value = 1
$ws.Cells.Item(114,9).Value = 0
$ws.Cells.Item(114,10).Value = "H"
$ws.Cells.Item(114,11).Value = 1.45
$ws.Cells.Item(114,12).Value = 4
$ws.Cells.Item(114,13).Value = 6
$ws.Cells.Item(114,14).Value = 1.333
$ws.Cells.Item(114,15).Value = 4.75
$ws.Cells.Item(114,16).Value = 7.5
$ws.Cells.Item(114,17).Value = -1.5
$ws.Cells.Item(114,18).Value = 1.95
$ws.Cells.Item(114,19).Value = 1.85
$ws.Cells.Item(114,20).Value = 3
$ws.Cells.Item(114,21).Value = 1.925
$ws.Cells.Item(114,22).Value = 1.875
$ws.Cells.Item(114,23).Value = 0.333
$ws.Cells.Item(114,24).Value = -1
$ws.Cells.Item(114,25).Value = -1
$ws.Cells.Item(114,26).Value = -1
$ws.Cells.Item(114,27).Value = 0.8500000000000001
$ws.Cells.Item(114,28).Value = -1
$ws.Cells.Item(114,29).Value = 0.875

# Row 115
$ws.Cells.Item(111,1).Copy() | Out-Null
$ws.Cells.Item(115,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(111,5).Copy() | Out-Null
$ws.Cells.Item(115,5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(115,1).Value = 113
$ws.Cells.Item(115,2).Value = 6859070
$ws.Cells.Item(115,3).Value = "Denmark Division 2"
$ws.Cells.Item(115,4).Value = "Denmark Division 2"
$ws.Cells.Item(115,5).Value = 45360.45833333334
$ws.Cells.Item(115,6).Value = "Skive"
$ws.Cells.Item(115,7).Value = "Fremad Amager"
$ws.Cells.Item(115,8).Value = 1
$ws.Cells.Item(115,9).Value = 1
$ws.Cells.Item(115,10).Value = "D"
$ws.Cells.Item(115,11).Value = 2.1
$ws.Cells.Item(115,12).Value = 3.6
$ws.Cells.Item(115,13).Value = 2.9
$ws.Cells.Item(115,14).Value = 1.85
$ws.Cells.Item(115,15).Value = 3.8
$ws.Cells.Item(115,16).Value = 3.4
$ws.Cells.Item(115,17).Value = -0.5
$ws.Cells.Item(115,18).Value = 1.875
$ws.Cells.Item(115,19).Value = 1.925
$ws.Cells.Item(115,20).Value = 2.25
$ws.Cells.Item(115,21).Value = 1.925
$ws.Cells.Item(115,22).Value = 1.875
$ws.Cells.Item(115,23).Value = -1
$ws.Cells.Item(115,24).Value = 2.8
$ws.Cells.Item(115,25).Value = -1
$ws.Cells.Item(115,26).Value = -1
$ws.Cells.Item(115,27).Value = 0.925
$ws.Cells.Item(115,28).Value = -0.5
$ws.Cells.Item(115,29).Value = 0.4375

# Row 116
$ws.Cells.Item(111,1).Copy() | Out-Null
$ws.Cells.Item(116,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(111,5).Copy() | Out-Null
$ws.Cells.Item(116,5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(116,1).Value = 114
$ws.Cells.Item(116,2).Value = 6859071
$ws.Cells.Item(116,3).Value = "Denmark Division 2"
$ws.Cells.Item(116,4).Value = "Denmark Division 2"
$ws.Cells.Item(116,5).Value = 45360.45833333334
$ws.Cells.Item(116,6).Value = "Nykobing"
$ws.Cells.Item(116,7).Value = "Hellerup IK"
$ws.Cells.Item(116,8).Value = 1
$ws.Cells.Item(116,9).Value = 0
$ws.Cells.Item(116,10).Value = "H"
$ws.Cells.Item(116,11).Value = 1.714
$ws.Cells.Item(116,12).Value = 3.75
$ws.Cells.Item(116,13).Value = 3.9
$ws.Cells.Item(116,14).Value = 1.85
$ws.Cells.Item(116,15).Value = 3.6
$ws.Cells.Item(116,16).Value = 3.4
$ws.Cells.Item(116,17).Value = -0.5
$ws.Cells.Item(116,18).Value = 1.9
$ws.Cells.Item(116,19).Value = 1.9
$ws.Cells.Item(116,20).Value = 2.75
$ws.Cells.Item(116,21).Value = 1.85
$ws.Cells.Item(116,22).Value = 1.95
$ws.Cells.Item(116,23).Value = 0.8500000000000001
$ws.Cells.Item(116,24).Value = -1
$ws.Cells.Item(116,25).Value = -1
$ws.Cells.Item(116,26).Value = 0.8999999999999999
$ws.Cells.Item(116,27).Value = -1
$ws.Cells.Item(116,28).Value = -1
$ws.Cells.Item(116,29).Value = 0.95

# Row 117
$ws.Cells.Item(111,1).Copy() | Out-Null
$ws.Cells.Item(117,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(111,5).Copy() | Out-Null
$ws.Cells.Item(117,5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(117,1).Value = 115
$ws.Cells.Item(117,2).Value = 6859069
$ws.Cells.Item(117,3).Value = "Denmark Division 2"
$ws.Cells.Item(117,4).Value = "Denmark Division 2"
$ws.Cells.Item(117,5).Value = 45361.41666666666
$ws.Cells.Item(117,6).Value = "Esbjerg"
$ws.Cells.Item(117,7).Value = "AB Copenhagen"
$ws.Cells.Item(117,8).Value = 2
$ws.Cells.Item(117,9).Value = 1
$ws.Cells.Item(117,10).Value = "H"
$ws.Cells.Item(117,11).Value = 1.25
$ws.Cells.Item(117,12).Value = 5.25
$ws.Cells.Item(117,13).Value = 11
$ws.Cells.Item(117,14).Value = 1.285
$ws.Cells.Item(117,15).Value = 5.25
$ws.Cells.Item(117,16).Value = 7.5
$ws.Cells.Item(117,17).Value = -1.5
$ws.Cells.Item(117,18).Value = 1.875
$ws.Cells.Item(117,19).Value = 1.925
$ws.Cells.Item(117,20).Value = 3
$ws.Cells.Item(117,21).Value = 1.85
$ws.Cells.Item(117,22).Value = 1.95
$ws.Cells.Item(117,23).Value = 0.2849999999999999
$ws.Cells.Item(117,24).Value = -1
$ws.Cells.Item(117,25).Value = -1
$ws.Cells.Item(117,26).Value = -1
$ws.Cells.Item(117,27).Value = 0.925
$ws.Cells.Item(117,28).Value = 0
$ws.Cells.Item(117,29).Value = -0

$excel.Application.CutCopyMode = 0